# Applies the "Atualizacao de bases das ligas" update to the Mexico Liga MX sheet.
# 1) 15 pairs of adjacent fixture rows swap their match data (sort-order change).
# 2) The fixture row that used to sit at row 410 (id 7612853) is removed; every
#    subsequent fixture row (411-418) shifts its data up by one row, and the row
#    that used to be the last one (418, id 7854021) is deleted outright.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31: swap with row 32
$ws.Range("B31").Value = 6001946
$ws.Range("F31").Value = "Club America"
$ws.Range("G31").Value = "Mazatlan FC"
$ws.Range("H31").Value = 6
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = "H"
$ws.Range("K31").Value = 1.3
$ws.Range("L31").Value = 5
$ws.Range("M31").Value = 10
$ws.Range("N31").Value = 1.222
$ws.Range("O31").Value = 6
$ws.Range("P31").Value = 12
$ws.Range("Q31").Value = -1.75
$ws.Range("R31").Value = 1.85
$ws.Range("S31").Value = 2
$ws.Range("T31").Value = 3.25
$ws.Range("U31").Value = 1.925
$ws.Range("V31").Value = 1.925
$ws.Range("W31").Value = 0.222
$ws.Range("Y31").Value = -1
$ws.Range("Z31").Value = 0.8500000000000001
$ws.Range("AA31").Value = -1
$ws.Range("AB31").Value = 0.925

# Row 32: swap with row 31
$ws.Range("B32").Value = 6001406
$ws.Range("F32").Value = "Juarez FC"
$ws.Range("G32").Value = "Chivas Guadalajara"
$ws.Range("H32").Value = 1
$ws.Range("I32").Value = 2
$ws.Range("J32").Value = "A"
$ws.Range("K32").Value = 2.5
$ws.Range("L32").Value = 3.3
$ws.Range("M32").Value = 2.75
$ws.Range("N32").Value = 2.875
$ws.Range("O32").Value = 2.875
$ws.Range("P32").Value = 2.75
$ws.Range("Q32").Value = 0
$ws.Range("R32").Value = 2.025
$ws.Range("S32").Value = 1.825
$ws.Range("T32").Value = 2
$ws.Range("U32").Value = 2
$ws.Range("V32").Value = 1.85
$ws.Range("W32").Value = -1
$ws.Range("Y32").Value = 1.75
$ws.Range("Z32").Value = -1
$ws.Range("AA32").Value = 0.825
$ws.Range("AB32").Value = 1

# Row 50: swap with row 51
$ws.Range("B50").Value = 6001964
$ws.Range("F50").Value = "Juarez FC"
$ws.Range("G50").Value = "Santos Laguna"
$ws.Range("H50").Value = 3
$ws.Range("I50").Value = 1
$ws.Range("K50").Value = 2.55
$ws.Range("L50").Value = 3.3
$ws.Range("M50").Value = 2.5
$ws.Range("N50").Value = 2.75
$ws.Range("O50").Value = 3.4
$ws.Range("P50").Value = 2.3
$ws.Range("Q50").Value = 0.25
$ws.Range("R50").Value = 1.8
$ws.Range("S50").Value = 2.05
$ws.Range("U50").Value = 2
$ws.Range("V50").Value = 1.85
$ws.Range("W50").Value = 1.75
$ws.Range("Z50").Value = 0.8
$ws.Range("AB50").Value = 1

# Row 51: swap with row 50
$ws.Range("B51").Value = 6001965
$ws.Range("F51").Value = "Tigres UANL"
$ws.Range("G51").Value = "Unam Pumas"
$ws.Range("H51").Value = 4
$ws.Range("I51").Value = 2
$ws.Range("K51").Value = 1.65
$ws.Range("L51").Value = 3.75
$ws.Range("M51").Value = 5
$ws.Range("N51").Value = 1.4
$ws.Range("O51").Value = 4.5
$ws.Range("P51").Value = 7
$ws.Range("Q51").Value = -1.25
$ws.Range("R51").Value = 1.9
$ws.Range("S51").Value = 1.95
$ws.Range("U51").Value = 1.925
$ws.Range("V51").Value = 1.925
$ws.Range("W51").Value = 0.3999999999999999
$ws.Range("Z51").Value = 0.8999999999999999
$ws.Range("AB51").Value = 0.925

# Row 54: swap with row 55
$ws.Range("B54").Value = 6001968
$ws.Range("F54").Value = "Tigres UANL"
$ws.Range("G54").Value = "Juarez FC"
$ws.Range("H54").Value = 0
$ws.Range("J54").Value = "D"
$ws.Range("K54").Value = 1.444
$ws.Range("L54").Value = 4.75
$ws.Range("M54").Value = 6
$ws.Range("N54").Value = 1.4
$ws.Range("O54").Value = 5
$ws.Range("P54").Value = 6.5
$ws.Range("Q54").Value = -1.25
$ws.Range("R54").Value = 1.975
$ws.Range("S54").Value = 1.875
$ws.Range("U54").Value = 1.825
$ws.Range("V54").Value = 2.025
$ws.Range("W54").Value = -1
$ws.Range("X54").Value = 4
$ws.Range("Z54").Value = -1
$ws.Range("AA54").Value = 0.875
$ws.Range("AC54").Value = 1.025

# Row 55: swap with row 54
$ws.Range("B55").Value = 6001967
$ws.Range("F55").Value = "Leon"
$ws.Range("G55").Value = "Puebla"
$ws.Range("H55").Value = 2
$ws.Range("J55").Value = "H"
$ws.Range("K55").Value = 2.05
$ws.Range("L55").Value = 3.4
$ws.Range("M55").Value = 3.6
$ws.Range("N55").Value = 1.65
$ws.Range("O55").Value = 4
$ws.Range("P55").Value = 5
$ws.Range("Q55").Value = -0.75
$ws.Range("R55").Value = 1.825
$ws.Range("S55").Value = 2.025
$ws.Range("U55").Value = 1.85
$ws.Range("V55").Value = 2
$ws.Range("W55").Value = 0.6499999999999999
$ws.Range("X55").Value = -1
$ws.Range("Z55").Value = 0.825
$ws.Range("AA55").Value = -1
$ws.Range("AC55").Value = 1

# Row 149: swap with row 150
$ws.Range("B149").Value = 6001419
$ws.Range("F149").Value = "Chivas Guadalajara"
$ws.Range("G149").Value = "Mazatlan FC"
$ws.Range("K149").Value = 1.363
$ws.Range("L149").Value = 4.75
$ws.Range("M149").Value = 7.5
$ws.Range("N149").Value = 1.333
$ws.Range("O149").Value = 5.25
$ws.Range("P149").Value = 9
$ws.Range("Q149").Value = -1.5
$ws.Range("R149").Value = 2
$ws.Range("S149").Value = 1.85
$ws.Range("T149").Value = 3
$ws.Range("U149").Value = 2.025
$ws.Range("V149").Value = 1.825
$ws.Range("W149").Value = 0.333
$ws.Range("Z149").Value = 1
$ws.Range("AB149").Value = 1.025

# Row 150: swap with row 149
$ws.Range("B150").Value = 6002050
$ws.Range("F150").Value = "Monterrey"
$ws.Range("G150").Value = "Unam Pumas"
$ws.Range("K150").Value = 1.55
$ws.Range("L150").Value = 4.333
$ws.Range("M150").Value = 5
$ws.Range("N150").Value = 1.7
$ws.Range("O150").Value = 4
$ws.Range("P150").Value = 4.75
$ws.Range("Q150").Value = -0.75
$ws.Range("R150").Value = 1.875
$ws.Range("S150").Value = 1.975
$ws.Range("T150").Value = 2.75
$ws.Range("U150").Value = 1.85
$ws.Range("V150").Value = 2
$ws.Range("W150").Value = 0.7
$ws.Range("Z150").Value = 0.875
$ws.Range("AB150").Value = 0.8500000000000001

# Row 205: swap with row 206
$ws.Range("B205").Value = 7053868
$ws.Range("F205").Value = "Necaxa"
$ws.Range("G205").Value = "Tigres UANL"
$ws.Range("H205").Value = 0
$ws.Range("I205").Value = 3
$ws.Range("J205").Value = "A"
$ws.Range("K205").Value = 3.75
$ws.Range("L205").Value = 3.3
$ws.Range("M205").Value = 1.909
$ws.Range("N205").Value = 4.5
$ws.Range("O205").Value = 3.6
$ws.Range("P205").Value = 1.833
$ws.Range("Q205").Value = 0.5
$ws.Range("R205").Value = 2.025
$ws.Range("S205").Value = 1.825
$ws.Range("T205").Value = 2.5
$ws.Range("X205").Value = -1
$ws.Range("Y205").Value = 0.833
$ws.Range("Z205").Value = -1
$ws.Range("AA205").Value = 0.825
$ws.Range("AB205").Value = 1
$ws.Range("AC205").Value = -1

# Row 206: swap with row 205
$ws.Range("B206").Value = 7094230
$ws.Range("F206").Value = "Club America"
$ws.Range("G206").Value = "Atlas"
$ws.Range("H206").Value = 1
$ws.Range("I206").Value = 1
$ws.Range("J206").Value = "D"
$ws.Range("K206").Value = 2.3
$ws.Range("L206").Value = 3.4
$ws.Range("M206").Value = 3.1
$ws.Range("N206").Value = 2.15
$ws.Range("O206").Value = 3.4
$ws.Range("P206").Value = 3.4
$ws.Range("Q206").Value = -0.25
$ws.Range("R206").Value = 1.85
$ws.Range("S206").Value = 2
$ws.Range("T206").Value = 2.75
$ws.Range("X206").Value = 2.4
$ws.Range("Y206").Value = -1
$ws.Range("Z206").Value = -0.5
$ws.Range("AA206").Value = 0.5
$ws.Range("AB206").Value = -1
$ws.Range("AC206").Value = 0.8500000000000001

# Row 243: swap with row 244
$ws.Range("B243").Value = 6754049
$ws.Range("F243").Value = "Juarez FC"
$ws.Range("G243").Value = "Atlas"
$ws.Range("H243").Value = 1
$ws.Range("J243").Value = "A"
$ws.Range("K243").Value = 2.75
$ws.Range("L243").Value = 3.25
$ws.Range("M243").Value = 2.375
$ws.Range("N243").Value = 2.6
$ws.Range("O243").Value = 3.2
$ws.Range("P243").Value = 2.8
$ws.Range("Q243").Value = 0
$ws.Range("R243").Value = 1.85
$ws.Range("S243").Value = 2
$ws.Range("T243").Value = 2.25
$ws.Range("U243").Value = 2.1
$ws.Range("V243").Value = 1.775
$ws.Range("W243").Value = -1
$ws.Range("Y243").Value = 1.8
$ws.Range("Z243").Value = -1
$ws.Range("AA243").Value = 1
$ws.Range("AB243").Value = 1.1

# Row 244: swap with row 243
$ws.Range("B244").Value = 6754048
$ws.Range("F244").Value = "Atletico San Luis"
$ws.Range("G244").Value = "Mazatlan FC"
$ws.Range("H244").Value = 3
$ws.Range("J244").Value = "H"
$ws.Range("K244").Value = 1.615
$ws.Range("L244").Value = 4
$ws.Range("M244").Value = 4.5
$ws.Range("N244").Value = 1.6
$ws.Range("O244").Value = 4.5
$ws.Range("P244").Value = 5
$ws.Range("Q244").Value = -1
$ws.Range("R244").Value = 1.95
$ws.Range("S244").Value = 1.9
$ws.Range("T244").Value = 3
$ws.Range("U244").Value = 1.925
$ws.Range("V244").Value = 1.925
$ws.Range("W244").Value = 0.6000000000000001
$ws.Range("Y244").Value = -1
$ws.Range("Z244").Value = 0
$ws.Range("AA244").Value = -0
$ws.Range("AB244").Value = 0.925

# Row 246: swap with row 247
$ws.Range("B246").Value = 6754051
$ws.Range("F246").Value = "Leon"
$ws.Range("G246").Value = "Tijuana"
$ws.Range("H246").Value = 1
$ws.Range("J246").Value = "H"
$ws.Range("K246").Value = 1.571
$ws.Range("L246").Value = 4
$ws.Range("M246").Value = 4.75
$ws.Range("N246").Value = 1.5
$ws.Range("O246").Value = 4.75
$ws.Range("P246").Value = 6
$ws.Range("Q246").Value = -1.25
$ws.Range("R246").Value = 2
$ws.Range("S246").Value = 1.85
$ws.Range("T246").Value = 3
$ws.Range("U246").Value = 1.925
$ws.Range("V246").Value = 1.925
$ws.Range("W246").Value = 0.5
$ws.Range("X246").Value = -1
$ws.Range("AA246").Value = 0.425
$ws.Range("AC246").Value = 0.925

# Row 247: swap with row 246
$ws.Range("B247").Value = 6754052
$ws.Range("F247").Value = "Chivas Guadalajara"
$ws.Range("G247").Value = "Pachuca"
$ws.Range("H247").Value = 0
$ws.Range("J247").Value = "D"
$ws.Range("K247").Value = 2
$ws.Range("L247").Value = 3.3
$ws.Range("M247").Value = 3.4
$ws.Range("N247").Value = 2.2
$ws.Range("O247").Value = 3.2
$ws.Range("P247").Value = 3.5
$ws.Range("Q247").Value = -0.25
$ws.Range("R247").Value = 1.875
$ws.Range("S247").Value = 1.975
$ws.Range("T247").Value = 2.5
$ws.Range("U247").Value = 2.025
$ws.Range("V247").Value = 1.825
$ws.Range("W247").Value = -1
$ws.Range("X247").Value = 2.2
$ws.Range("AA247").Value = 0.4875
$ws.Range("AC247").Value = 0.825

# Row 263: swap with row 264
$ws.Range("B263").Value = 6754066
$ws.Range("F263").Value = "Unam Pumas"
$ws.Range("G263").Value = "Queretaro"
$ws.Range("H263").Value = 4
$ws.Range("I263").Value = 0
$ws.Range("J263").Value = "H"
$ws.Range("K263").Value = 1.727
$ws.Range("L263").Value = 3.5
$ws.Range("M263").Value = 4.5
$ws.Range("N263").Value = 1.8
$ws.Range("P263").Value = 4.5
$ws.Range("Q263").Value = -0.75
$ws.Range("R263").Value = 2.025
$ws.Range("S263").Value = 1.825
$ws.Range("U263").Value = 1.825
$ws.Range("V263").Value = 2.025
$ws.Range("W263").Value = 0.8
$ws.Range("Y263").Value = -1
$ws.Range("Z263").Value = 1.025
$ws.Range("AA263").Value = -1
$ws.Range("AB263").Value = 0.825

# Row 264: swap with row 263
$ws.Range("B264").Value = 6754065
$ws.Range("F264").Value = "Necaxa"
$ws.Range("G264").Value = "Cruz Azul"
$ws.Range("H264").Value = 1
$ws.Range("I264").Value = 3
$ws.Range("J264").Value = "A"
$ws.Range("K264").Value = 2.375
$ws.Range("L264").Value = 3.3
$ws.Range("M264").Value = 2.8
$ws.Range("N264").Value = 3.5
$ws.Range("P264").Value = 2.1
$ws.Range("Q264").Value = 0.25
$ws.Range("R264").Value = 2
$ws.Range("S264").Value = 1.85
$ws.Range("U264").Value = 1.9
$ws.Range("V264").Value = 1.95
$ws.Range("W264").Value = -1
$ws.Range("Y264").Value = 1.1
$ws.Range("Z264").Value = -1
$ws.Range("AA264").Value = 0.8500000000000001
$ws.Range("AB264").Value = 0.8999999999999999

# Row 265: swap with row 266
$ws.Range("B265").Value = 7260442
$ws.Range("F265").Value = "Santos Laguna"
$ws.Range("G265").Value = "Tijuana"
$ws.Range("I265").Value = 1
$ws.Range("J265").Value = "H"
$ws.Range("K265").Value = 1.75
$ws.Range("L265").Value = 3.6
$ws.Range("M265").Value = 4.2
$ws.Range("N265").Value = 1.65
$ws.Range("O265").Value = 4
$ws.Range("P265").Value = 4.75
$ws.Range("Q265").Value = -0.75
$ws.Range("R265").Value = 1.8
$ws.Range("S265").Value = 2.05
$ws.Range("U265").Value = 1.85
$ws.Range("V265").Value = 2
$ws.Range("W265").Value = 0.6499999999999999
$ws.Range("X265").Value = -1
$ws.Range("Z265").Value = 0.4
$ws.Range("AA265").Value = -0.5
$ws.Range("AB265").Value = 0
$ws.Range("AC265").Value = -0

# Row 266: swap with row 265
$ws.Range("B266").Value = 6754067
$ws.Range("F266").Value = "Tigres UANL"
$ws.Range("G266").Value = "Toluca"
$ws.Range("I266").Value = 2
$ws.Range("J266").Value = "D"
$ws.Range("K266").Value = 1.8
$ws.Range("L266").Value = 3.3
$ws.Range("M266").Value = 4.333
$ws.Range("N266").Value = 1.533
$ws.Range("O266").Value = 4.2
$ws.Range("P266").Value = 6
$ws.Range("Q266").Value = -1
$ws.Range("R266").Value = 1.925
$ws.Range("S266").Value = 1.925
$ws.Range("U266").Value = 2
$ws.Range("V266").Value = 1.85
$ws.Range("W266").Value = -1
$ws.Range("X266").Value = 3.2
$ws.Range("Z266").Value = -1
$ws.Range("AA266").Value = 0.925
$ws.Range("AB266").Value = 1
$ws.Range("AC266").Value = -1

# Row 269: swap with row 270
$ws.Range("B269").Value = 6754074
$ws.Range("F269").Value = "Chivas Guadalajara"
$ws.Range("G269").Value = "Atlas"
$ws.Range("H269").Value = 4
$ws.Range("J269").Value = "H"
$ws.Range("K269").Value = 2.3
$ws.Range("L269").Value = 3.3
$ws.Range("M269").Value = 2.8
$ws.Range("N269").Value = 2.4
$ws.Range("O269").Value = 3.2
$ws.Range("P269").Value = 3.1
$ws.Range("Q269").Value = -0.25
$ws.Range("R269").Value = 2.15
$ws.Range("S269").Value = 1.725
$ws.Range("T269").Value = 2
$ws.Range("U269").Value = 1.925
$ws.Range("V269").Value = 1.925
$ws.Range("W269").Value = 1.4
$ws.Range("X269").Value = -1
$ws.Range("Z269").Value = 1.15
$ws.Range("AA269").Value = -1
$ws.Range("AB269").Value = 0.925
$ws.Range("AC269").Value = -1

# Row 270: swap with row 269
$ws.Range("B270").Value = 6754641
$ws.Range("F270").Value = "Pachuca"
$ws.Range("G270").Value = "Tigres UANL"
$ws.Range("H270").Value = 1
$ws.Range("J270").Value = "D"
$ws.Range("K270").Value = 2.875
$ws.Range("L270").Value = 3.5
$ws.Range("M270").Value = 2.2
$ws.Range("N270").Value = 2.9
$ws.Range("O270").Value = 3.5
$ws.Range("P270").Value = 2.4
$ws.Range("Q270").Value = 0.25
$ws.Range("R270").Value = 1.775
$ws.Range("S270").Value = 2.1
$ws.Range("T270").Value = 2.5
$ws.Range("U270").Value = 1.825
$ws.Range("V270").Value = 2.025
$ws.Range("W270").Value = -1
$ws.Range("X270").Value = 2.5
$ws.Range("Z270").Value = 0.3875
$ws.Range("AA270").Value = -0.5
$ws.Range("AB270").Value = -1
$ws.Range("AC270").Value = 1.025

# Row 297: swap with row 298
$ws.Range("B297").Value = 6754096
$ws.Range("F297").Value = "Queretaro"
$ws.Range("G297").Value = "Chivas Guadalajara"
$ws.Range("H297").Value = 1
$ws.Range("I297").Value = 2
$ws.Range("K297").Value = 3.25
$ws.Range("L297").Value = 3.3
$ws.Range("M297").Value = 2.2
$ws.Range("N297").Value = 2.8
$ws.Range("O297").Value = 3.1
$ws.Range("P297").Value = 2.7
$ws.Range("Q297").Value = 0
$ws.Range("T297").Value = 2.25
$ws.Range("U297").Value = 1.925
$ws.Range("V297").Value = 1.925
$ws.Range("Y297").Value = 1.7
$ws.Range("AB297").Value = 0.925
$ws.Range("AC297").Value = -1

# Row 298: swap with row 297
$ws.Range("B298").Value = 6754097
$ws.Range("F298").Value = "Toluca"
$ws.Range("G298").Value = "Puebla"
$ws.Range("H298").Value = 0
$ws.Range("I298").Value = 1
$ws.Range("K298").Value = 1.5
$ws.Range("L298").Value = 4
$ws.Range("M298").Value = 7
$ws.Range("N298").Value = 1.45
$ws.Range("O298").Value = 4.5
$ws.Range("P298").Value = 7
$ws.Range("Q298").Value = -1.25
$ws.Range("T298").Value = 3.25
$ws.Range("U298").Value = 1.975
$ws.Range("V298").Value = 1.875
$ws.Range("Y298").Value = 6
$ws.Range("AB298").Value = -1
$ws.Range("AC298").Value = 0.875

# Row 301: swap with row 302
$ws.Range("B301").Value = 6754101
$ws.Range("F301").Value = "Cruz Azul"
$ws.Range("G301").Value = "Juarez FC"
$ws.Range("H301").Value = 2
$ws.Range("I301").Value = 0
$ws.Range("J301").Value = "H"
$ws.Range("K301").Value = 2.25
$ws.Range("M301").Value = 3.25
$ws.Range("N301").Value = 1.65
$ws.Range("O301").Value = 3.8
$ws.Range("P301").Value = 5.25
$ws.Range("Q301").Value = -0.75
$ws.Range("R301").Value = 1.8
$ws.Range("S301").Value = 2.05
$ws.Range("T301").Value = 2.75
$ws.Range("U301").Value = 1.925
$ws.Range("V301").Value = 1.925
$ws.Range("W301").Value = 0.6499999999999999
$ws.Range("Y301").Value = -1
$ws.Range("Z301").Value = 0.8
$ws.Range("AA301").Value = -1
$ws.Range("AB301").Value = -1
$ws.Range("AC301").Value = 0.925

# Row 302: swap with row 301
$ws.Range("B302").Value = 6754100
$ws.Range("F302").Value = "Atlas"
$ws.Range("G302").Value = "Pachuca"
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 2
$ws.Range("J302").Value = "A"
$ws.Range("K302").Value = 2.1
$ws.Range("M302").Value = 3.5
$ws.Range("N302").Value = 2.45
$ws.Range("O302").Value = 3.2
$ws.Range("P302").Value = 3
$ws.Range("Q302").Value = -0.25
$ws.Range("R302").Value = 2.05
$ws.Range("S302").Value = 1.8
$ws.Range("T302").Value = 2.25
$ws.Range("U302").Value = 1.875
$ws.Range("V302").Value = 1.975
$ws.Range("W302").Value = -1
$ws.Range("Y302").Value = 2
$ws.Range("Z302").Value = -1
$ws.Range("AA302").Value = 0.8
$ws.Range("AB302").Value = -0.5
$ws.Range("AC302").Value = 0.4875

# Row 371: swap with row 372
$ws.Range("B371").Value = 7612821
$ws.Range("F371").Value = "Mazatlan FC"
$ws.Range("G371").Value = "Leon"
$ws.Range("H371").Value = 2
$ws.Range("I371").Value = 2
$ws.Range("J371").Value = "D"
$ws.Range("K371").Value = 2.5
$ws.Range("L371").Value = 3.3
$ws.Range("M371").Value = 2.75
$ws.Range("N371").Value = 3.3
$ws.Range("O371").Value = 3.6
$ws.Range("P371").Value = 2.15
$ws.Range("Q371").Value = 0.25
$ws.Range("R371").Value = 1.925
$ws.Range("S371").Value = 1.925
$ws.Range("T371").Value = 2.5
$ws.Range("U371").Value = 1.875
$ws.Range("V371").Value = 1.975
$ws.Range("W371").Value = -1
$ws.Range("X371").Value = 2.6
$ws.Range("Z371").Value = 0.4625
$ws.Range("AA371").Value = -0.5
$ws.Range("AB371").Value = 0.875
$ws.Range("AC371").Value = -1

# Row 372: swap with row 371
$ws.Range("B372").Value = 7713694
$ws.Range("F372").Value = "Cruz Azul"
$ws.Range("G372").Value = "Tijuana"
$ws.Range("H372").Value = 1
$ws.Range("I372").Value = 0
$ws.Range("J372").Value = "H"
$ws.Range("K372").Value = 1.727
$ws.Range("L372").Value = 3.75
$ws.Range("M372").Value = 4.5
$ws.Range("N372").Value = 1.533
$ws.Range("O372").Value = 4.2
$ws.Range("P372").Value = 6
$ws.Range("Q372").Value = -1
$ws.Range("R372").Value = 1.95
$ws.Range("S372").Value = 1.9
$ws.Range("T372").Value = 2.75
$ws.Range("U372").Value = 1.95
$ws.Range("V372").Value = 1.9
$ws.Range("W372").Value = 0.5329999999999999
$ws.Range("X372").Value = -1
$ws.Range("Z372").Value = 0
$ws.Range("AA372").Value = -0
$ws.Range("AB372").Value = -1
$ws.Range("AC372").Value = 0.8999999999999999

# Row 393: swap with row 394
$ws.Range("B393").Value = 7612842
$ws.Range("F393").Value = "Santos Laguna"
$ws.Range("G393").Value = "Tigres UANL"
$ws.Range("H393").Value = 0
$ws.Range("I393").Value = 3
$ws.Range("J393").Value = "A"
$ws.Range("K393").Value = 3
$ws.Range("L393").Value = 3.5
$ws.Range("M393").Value = 2.15
$ws.Range("N393").Value = 3.4
$ws.Range("O393").Value = 3.5
$ws.Range("P393").Value = 2.1
$ws.Range("Q393").Value = 0.25
$ws.Range("R393").Value = 2.05
$ws.Range("S393").Value = 1.8
$ws.Range("T393").Value = 2.5
$ws.Range("U393").Value = 1.9
$ws.Range("V393").Value = 1.95
$ws.Range("W393").Value = -1
$ws.Range("Y393").Value = 1.1
$ws.Range("Z393").Value = -1
$ws.Range("AA393").Value = 0.8
$ws.Range("AB393").Value = 0.8999999999999999

# Row 394: swap with row 393
$ws.Range("B394").Value = 7612841
$ws.Range("F394").Value = "Monterrey"
$ws.Range("G394").Value = "Pachuca"
$ws.Range("H394").Value = 3
$ws.Range("I394").Value = 2
$ws.Range("J394").Value = "H"
$ws.Range("K394").Value = 1.5
$ws.Range("L394").Value = 4.2
$ws.Range("M394").Value = 5.5
$ws.Range("N394").Value = 1.7
$ws.Range("O394").Value = 4
$ws.Range("P394").Value = 4.75
$ws.Range("Q394").Value = -0.75
$ws.Range("R394").Value = 1.95
$ws.Range("S394").Value = 1.9
$ws.Range("T394").Value = 2.75
$ws.Range("U394").Value = 1.825
$ws.Range("V394").Value = 2.025
$ws.Range("W394").Value = 0.7
$ws.Range("Y394").Value = -1
$ws.Range("Z394").Value = 0.475
$ws.Range("AA394").Value = -0.5
$ws.Range("AB394").Value = 0.825

# Row 408: swap with row 409
$ws.Range("B408").Value = 7612867
$ws.Range("F408").Value = "Club America"
$ws.Range("G408").Value = "Mazatlan FC"
$ws.Range("I408").Value = 2
$ws.Range("J408").Value = "D"
$ws.Range("K408").Value = 1.363
$ws.Range("L408").Value = 5
$ws.Range("M408").Value = 7.5
$ws.Range("N408").Value = 1.222
$ws.Range("O408").Value = 6.5
$ws.Range("P408").Value = 12
$ws.Range("Q408").Value = -1.75
$ws.Range("R408").Value = 1.825
$ws.Range("S408").Value = 2.025
$ws.Range("T408").Value = 3.25
$ws.Range("U408").Value = 1.975
$ws.Range("V408").Value = 1.875
$ws.Range("X408").Value = 5.5
$ws.Range("Y408").Value = -1
$ws.Range("AA408").Value = 1.025
$ws.Range("AB408").Value = 0.9750000000000001

# Row 409: swap with row 408
$ws.Range("B409").Value = 7612866
$ws.Range("F409").Value = "Leon"
$ws.Range("G409").Value = "Cruz Azul"
$ws.Range("I409").Value = 3
$ws.Range("J409").Value = "A"
$ws.Range("K409").Value = 2.5
$ws.Range("L409").Value = 3.4
$ws.Range("M409").Value = 2.7
$ws.Range("N409").Value = 2.8
$ws.Range("O409").Value = 3.6
$ws.Range("P409").Value = 2.375
$ws.Range("Q409").Value = 0.25
$ws.Range("R409").Value = 1.75
$ws.Range("S409").Value = 2.05
$ws.Range("T409").Value = 2.75
$ws.Range("U409").Value = 1.85
$ws.Range("V409").Value = 2
$ws.Range("X409").Value = -1
$ws.Range("Y409").Value = 1.375
$ws.Range("AA409").Value = 1.05
$ws.Range("AB409").Value = 0.8500000000000001

# Row 410: take values from old row 411 (rows shift up after the row-410 deletion)
$ws.Range("B410").Value = 7612854
$ws.Range("E410").Value = 45346
$ws.Range("F410").Value = "Necaxa"
$ws.Range("G410").Value = "Pachuca"
$ws.Range("K410").Value = 2.75
$ws.Range("L410").Value = 3.3
$ws.Range("M410").Value = 2.45
$ws.Range("N410").Value = 3.2
$ws.Range("O410").Value = 3.5
$ws.Range("P410").Value = 2.2
$ws.Range("Q410").Value = 0.25

# Row 411: take values from old row 412 (rows shift up after the row-410 deletion)
$ws.Range("B411").Value = 7612855
$ws.Range("E411").Value = 45346.00694444445
$ws.Range("F411").Value = "Juarez FC"
$ws.Range("G411").Value = "Monterrey"
$ws.Range("K411").Value = 5
$ws.Range("L411").Value = 4.1
$ws.Range("M411").Value = 1.533
$ws.Range("N411").Value = 4.5
$ws.Range("O411").Value = 3.8
$ws.Range("P411").Value = 1.75
$ws.Range("Q411").Value = 0.75
$ws.Range("T411").Value = 2.5
$ws.Range("U411").Value = 1.925
$ws.Range("V411").Value = 1.925

# Row 412: take values from old row 413 (rows shift up after the row-410 deletion)
$ws.Range("B412").Value = 7612856
$ws.Range("E412").Value = 45346.83333333334
$ws.Range("F412").Value = "Leon"
$ws.Range("G412").Value = "Atletico San Luis"
$ws.Range("K412").Value = 2.1
$ws.Range("L412").Value = 3.5
$ws.Range("M412").Value = 3
$ws.Range("N412").Value = 1.85
$ws.Range("O412").Value = 4
$ws.Range("P412").Value = 3.8
$ws.Range("Q412").Value = -0.5
$ws.Range("R412").Value = 1.875
$ws.Range("S412").Value = 1.975
$ws.Range("T412").Value = 3
$ws.Range("U412").Value = 1.975
$ws.Range("V412").Value = 1.875

# Row 413: take values from old row 414 (rows shift up after the row-410 deletion)
$ws.Range("B413").Value = 7612858
$ws.Range("E413").Value = 45346.91666666666
$ws.Range("F413").Value = "Tigres UANL"
$ws.Range("G413").Value = "Atlas"
$ws.Range("K413").Value = 1.444
$ws.Range("L413").Value = 4.2
$ws.Range("M413").Value = 6.2
$ws.Range("N413").Value = 1.571
$ws.Range("O413").Value = 4.2
$ws.Range("P413").Value = 5.5
$ws.Range("Q413").Value = -1
$ws.Range("R413").Value = 2.025
$ws.Range("S413").Value = 1.825
$ws.Range("T413").Value = 2.5
$ws.Range("U413").Value = 2
$ws.Range("V413").Value = 1.85

# Row 414: take values from old row 415 (rows shift up after the row-410 deletion)
$ws.Range("B414").Value = 7612857
$ws.Range("E414").Value = 45346.92013888889
$ws.Range("F414").Value = "Chivas Guadalajara"
$ws.Range("G414").Value = "Unam Pumas"
$ws.Range("K414").Value = 2.3
$ws.Range("L414").Value = 3.4
$ws.Range("M414").Value = 2.75
$ws.Range("N414").Value = 2.15
$ws.Range("O414").Value = 3.4
$ws.Range("P414").Value = 3.3
$ws.Range("Q414").Value = -0.25
$ws.Range("R414").Value = 1.875
$ws.Range("S414").Value = 1.975
$ws.Range("U414").Value = 1.975
$ws.Range("V414").Value = 1.875

# Row 415: take values from old row 416 (rows shift up after the row-410 deletion)
$ws.Range("B415").Value = 7612859
$ws.Range("E415").Value = 45347
$ws.Range("F415").Value = "Club America"
$ws.Range("G415").Value = "Cruz Azul"
$ws.Range("K415").Value = 1.444
$ws.Range("L415").Value = 4.25
$ws.Range("M415").Value = 6
$ws.Range("N415").Value = 1.95
$ws.Range("O415").Value = 3.8
$ws.Range("P415").Value = 3.6
$ws.Range("Q415").Value = -0.5
$ws.Range("R415").Value = 1.95
$ws.Range("S415").Value = 1.9
$ws.Range("T415").Value = 2.75
$ws.Range("U415").Value = 1.95
$ws.Range("V415").Value = 1.9

# Row 416: take values from old row 417 (rows shift up after the row-410 deletion)
$ws.Range("B416").Value = 7612860
$ws.Range("E416").Value = 45347.625
$ws.Range("F416").Value = "Toluca"
$ws.Range("G416").Value = "Tijuana"
$ws.Range("K416").Value = 1.75
$ws.Range("L416").Value = 3.6
$ws.Range("M416").Value = 4.1
$ws.Range("N416").Value = 1.615
$ws.Range("O416").Value = 4.2
$ws.Range("P416").Value = 5
$ws.Range("Q416").Value = -1
$ws.Range("R416").Value = 2.05
$ws.Range("S416").Value = 1.8
$ws.Range("T416").Value = 3
$ws.Range("U416").Value = 1.925
$ws.Range("V416").Value = 1.925

# Row 417: take values from old row 418 (rows shift up after the row-410 deletion)
$ws.Range("B417").Value = 7854021
$ws.Range("E417").Value = 45347.875
$ws.Range("F417").Value = "Santos Laguna"
$ws.Range("G417").Value = "Mazatlan FC"
$ws.Range("K417").Value = 1.95
$ws.Range("M417").Value = 3.3
$ws.Range("N417").Value = 1.95
$ws.Range("O417").Value = 3.8
$ws.Range("P417").Value = 3.6
$ws.Range("Q417").Value = -0.5
$ws.Range("R417").Value = 1.975
$ws.Range("S417").Value = 1.875
$ws.Range("U417").Value = 2
$ws.Range("V417").Value = 1.85

# The old last row (418, id 7854021) no longer exists; its data already lives in
# row 417 above, so just remove the now-duplicate trailing row.
$ws.Rows.Item(418).Delete()
